$d = $word.ActiveDocument

$replacements = @(
    @{old = "547÷2="; new = "699÷7="},
    @{old = "814÷6="; new = "582÷3="},
    @{old = "982÷9="; new = "527÷2="},
    @{old = "929÷6="; new = "996÷5="},
    @{old = "306÷3="; new = "121÷8="},
    @{old = "202÷8="; new = "707÷3="},
    @{old = "260÷8="; new = "611÷2="},
    @{old = "125÷7="; new = "990÷9="},
    @{old = "981÷6="; new = "488÷9="},
    @{old = "685÷2="; new = "986÷8="},
    @{old = "478÷5="; new = "308÷5="},
    @{old = "105÷5="; new = "312÷3="},
    @{old = "336÷4="; new = "838÷9="},
    @{old = "557÷8="; new = "902÷5="},
    @{old = "914÷7="; new = "145÷8="},
    @{old = "776÷8="; new = "905÷5="},
    @{old = "773÷5="; new = "628÷2="},
    @{old = "953÷9="; new = "880÷7="},
    @{old = "212÷3="; new = "589÷2="},
    @{old = "247÷6="; new = "947÷2="},
    @{old = "107÷7="; new = "428÷6="},
    @{old = "868÷8="; new = "508÷6="},
    @{old = "447÷5="; new = "529÷6="},
    @{old = "153÷2="; new = "398÷5="},
    @{old = "139÷6="; new = "295÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
